$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J70").Value = 2956.3333
$ws.Range("L70").Value = 8868.999899999999
$ws.Range("H70").Value = 1744741.4
$ws.Range("N70").Value = -9408.999899999999
$ws.Range("L73").Value = 8868.999899999999
$ws.Range("J73").Value = 2956.3333
$ws.Range("H73").Value = 1744741.4
$ws.Range("N73").Value = -10740.9999
$ws.Range("J132").Value = 4199.25
$ws.Range("L132").Value = 12597.75
$ws.Range("K132").Value = 5838.4998
$ws.Range("I132").Value = 1946.1666
$ws.Range("H132").Value = 2355.818
$ws.Range("M132").Value = -3308.4998
$ws.Range("N132").Value = -17657.75
$ws.Range("H137").Value = 19233056
$ws.Range("K137").Value = 93755226
$ws.Range("M137").Value = -93752676
$ws.Range("I137").Value = 31251742
$ws.Range("H138").Value = 3267.3035
$ws.Range("M138").Value = 798.4546
$ws.Range("I138").Value = 1447.1818
$ws.Range("J138").Value = 4445.0293
$ws.Range("L138").Value = 13335.0879
$ws.Range("K138").Value = 4341.5454
$ws.Range("N138").Value = -23615.0879

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K45").Value = 1062.6
$ws.Range("H45").Value = 2190.95
$ws.Range("I45").Value = 1062.6
$ws.Range("M45").Value = -685.5999999999999
$ws.Range("H74").Value = 2257.6
$ws.Range("K74").Value = 2119.5557
$ws.Range("J74").Value = 3500
$ws.Range("M74").Value = -1245.5557
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5248
$ws.Range("I74").Value = 2119.5557
$ws.Range("N77").Value = -26236
$ws.Range("I77").Value = 2119.5557
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -6229.7785
$ws.Range("K77").Value = 10597.7785
$ws.Range("H77").Value = 2257.6
$ws.Range("J77").Value = 3500
$ws.Range("H139").Value = 140000
$ws.Range("L139").Value = 140000
$ws.Range("N139").Value = -150280
$ws.Range("J139").Value = 140000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K20").Value = 1977.7
$ws.Range("I20").Value = 1977.7
$ws.Range("H20").Value = 2264.1162
$ws.Range("M20").Value = -1730.7
$ws.Range("L94").Value = 1168.6
$ws.Range("J94").Value = 1168.6
$ws.Range("I94").Value = 1988.8148
$ws.Range("K94").Value = 1988.8148
$ws.Range("H94").Value = 1767.1351
$ws.Range("M94").Value = -1537.8148
$ws.Range("N94").Value = -2070.6
$ws.Range("K99").Value = 3224
$ws.Range("M99").Value = -1726
$ws.Range("H99").Value = 3432.889
$ws.Range("I99").Value = 3224
$ws.Range("I107").Value = 4961.1113
$ws.Range("M107").Value = -3041.1113
$ws.Range("K107").Value = 4961.1113
$ws.Range("H107").Value = 4715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50
$ws.Range("H22").Value = 400
$ws.Range("M31").Value = -25001673
$ws.Range("I31").Value = 25001968
$ws.Range("H31").Value = 18870416
$ws.Range("K31").Value = 25001968
$ws.Range("I34").Value = 25001968
$ws.Range("K34").Value = 25001968
$ws.Range("H34").Value = 18870416
$ws.Range("M34").Value = -25001766
$ws.Range("I58").Value = 2222.5833
$ws.Range("H58").Value = 2977.9333
$ws.Range("N58").Value = -6405.3335
$ws.Range("J58").Value = 5999.3335
$ws.Range("L58").Value = 5999.3335
$ws.Range("M58").Value = -2019.5833
$ws.Range("K58").Value = 2222.5833
$ws.Range("N86").Value = -11143.875
$ws.Range("J86").Value = 8897.875
$ws.Range("M86").Value = -5471.3
$ws.Range("H86").Value = 7618.1113
$ws.Range("L86").Value = 8897.875
$ws.Range("I86").Value = 6594.3
$ws.Range("K86").Value = 6594.3
$ws.Range("K89").Value = 32971.5
$ws.Range("N89").Value = -55721.375
$ws.Range("L89").Value = 44489.375
$ws.Range("J89").Value = 8897.875
$ws.Range("H89").Value = 7618.1113
$ws.Range("M89").Value = -27355.5
$ws.Range("I89").Value = 6594.3
$ws.Range("L94").Value = 1576.5
$ws.Range("J94").Value = 1576.5
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("H94").Value = 1576.5
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2478.5
$ws.Range("M105").Value = 469
$ws.Range("K105").Value = 1278
$ws.Range("I105").Value = 1278
$ws.Range("H105").Value = 9457.888999999999
$ws.Range("H134").Value = 1997.6285
$ws.Range("K134").Value = 4937.6538
$ws.Range("I134").Value = 1645.8846
$ws.Range("M134").Value = -2402.6538
$ws.Range("I136").Value = 2222.5833
$ws.Range("L136").Value = 17998.0005
$ws.Range("J136").Value = 5999.3335
$ws.Range("K136").Value = 6667.749899999999
$ws.Range("N136").Value = -23098.0005
$ws.Range("M136").Value = -4117.749899999999
$ws.Range("H136").Value = 2977.9333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 92840
$ws.Range("J37").Value = 92840
$ws.Range("N37").Value = -278744
$ws.Range("L37").Value = 278520
$ws.Range("H124").Value = 9530.166999999999
$ws.Range("I124").Value = 1925
$ws.Range("N124").Value = -49818.25
$ws.Range("M124").Value = -865
$ws.Range("K124").Value = 5775
$ws.Range("L124").Value = 39998.25
$ws.Range("J124").Value = 13332.75
$ws.Range("K132").Value = 22718.25
$ws.Range("I132").Value = 2524.25
$ws.Range("H132").Value = 2524.25
$ws.Range("M132").Value = -20188.25
$ws.Range("H137").Value = 19958.25
$ws.Range("N137").Value = -109699.5
$ws.Range("L137").Value = 99499.5
$ws.Range("J137").Value = 33166.5
$ws.Range("K139").Value = 12772.7505
$ws.Range("I139").Value = 4257.5835
$ws.Range("M139").Value = -7632.750499999998
$ws.Range("H139").Value = 6681.2666
$ws.Range("L139").Value = 49128
$ws.Range("N139").Value = -59408
$ws.Range("J139").Value = 16376

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L45").Value = 98217
$ws.Range("H45").Value = 98217
$ws.Range("N45").Value = -99335
$ws.Range("J45").Value = 98217
$ws.Range("K70").Value = 12021.6875
$ws.Range("H70").Value = 13750.297
$ws.Range("M70").Value = -11751.6875
$ws.Range("I70").Value = 12021.6875
$ws.Range("K73").Value = 12021.6875
$ws.Range("M73").Value = -11085.6875
$ws.Range("H73").Value = 13750.297
$ws.Range("I73").Value = 12021.6875
$ws.Range("N88").ClearContents()
$ws.Range("J88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("L94").Value = 62487
$ws.Range("J94").Value = 62487
$ws.Range("H94").Value = 62487
$ws.Range("N94").Value = -63839
$ws.Range("K132").Value = 10482.717
$ws.Range("I132").Value = 3494.239
$ws.Range("H132").Value = 3239006
$ws.Range("M132").Value = -7952.717000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 7926.706
$ws.Range("K7").Value = 7926.706
$ws.Range("H7").Value = 7695.1035
$ws.Range("M7").Value = -7814.706
$ws.Range("I16").Value = 2112.2104
$ws.Range("H16").Value = 2151.4583
$ws.Range("M16").Value = -1942.2104
$ws.Range("K16").Value = 2112.2104
$ws.Range("M100").Value = -1688.8572
$ws.Range("K100").Value = 2229.8572
$ws.Range("I100").Value = 2229.8572
$ws.Range("H100").Value = 13905620
$ws.Range("H126").Value = 7695.1035
$ws.Range("M126").Value = -21310.118
$ws.Range("I126").Value = 7926.706
$ws.Range("K126").Value = 23780.118
$ws.Range("K132").Value = 5853.8181
$ws.Range("I132").Value = 1951.2727
$ws.Range("H132").Value = 3439.8572
$ws.Range("M132").Value = -3323.8181
